# Edit LocalizedSettings sheet: reorder form-option rows before the Submit
# button, add a FormCancelButton row, add a StoppedExecution message row,
# resize Table13, and update the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LocalizedSettings")

$data = New-Object 'object[,]' 83,3
$data[0,0] = 'Name'
$data[0,1] = 'EN'
$data[0,2] = 'JA'
$data[1,0] = 'FormTitle'
$data[1,1] = 'Orchestrator Manager'
$data[1,2] = 'Orchestrator Manager'
$data[2,0] = 'FormUsernameLabel'
$data[2,1] = 'Username'
$data[2,2] = 'ユーザー名'
$data[3,0] = 'FormPasswordLabel'
$data[3,1] = 'Password'
$data[3,2] = 'パスワード'
$data[4,0] = 'FormEntityLabel'
$data[4,1] = 'Entity'
$data[4,2] = 'エンティティ'
$data[5,0] = 'FormOperationLabel'
$data[5,1] = 'Operation'
$data[5,2] = '操作'
$data[6,0] = 'FormAssetOption'
$data[6,1] = 'Asset'
$data[6,2] = 'アセット'
$data[7,0] = 'FormUserOption'
$data[7,1] = 'User'
$data[7,2] = 'ユーザー'
$data[8,0] = 'FormMachineOption'
$data[8,1] = 'Machine'
$data[8,2] = 'マシン'
$data[9,0] = 'FormRobotOption'
$data[9,1] = 'Robot'
$data[9,2] = 'ロボット'
$data[10,0] = 'FormEnvironmentOption'
$data[10,1] = 'Environment'
$data[10,2] = 'ロボットグループ'
$data[11,0] = 'FormSubmitButton'
$data[11,1] = 'Submit'
$data[11,2] = '送信'
$data[12,0] = 'FormCancelButton'
$data[12,1] = 'Cancel'
$data[12,2] = 'キャンセル'
$data[13,0] = ''
$data[13,1] = ''
$data[13,2] = ''
$data[14,0] = 'AssetConfigFilepath'
$data[14,1] = 'Config\EN\Assets.xlsx'
$data[14,2] = 'Config\JA\アセット.xlsx'
$data[15,0] = 'UserConfigFilepath'
$data[15,1] = 'Config\EN\Users.xlsx'
$data[15,2] = 'Config\JA\ユーザー.xlsx'
$data[16,0] = 'RobotConfigFilepath'
$data[16,1] = 'Config\EN\Robots.xlsx'
$data[16,2] = 'Config\JA\ロボット.xlsx'
$data[17,0] = 'MachineConfigFilepath'
$data[17,1] = 'Config\EN\Machines.xlsx'
$data[17,2] = 'Config\JA\マシン.xlsx'
$data[18,0] = 'EnvironmentConfigFilepath'
$data[18,1] = 'Config\EN\Environments.xlsx'
$data[18,2] = 'Config\JA\ロボットグループ.xlsx'
$data[19,0] = ''
$data[19,1] = ''
$data[19,2] = ''
$data[20,0] = 'GetOperationName'
$data[20,1] = 'Get'
$data[20,2] = '取得'
$data[21,0] = 'CreateOperationName'
$data[21,1] = 'Create'
$data[21,2] = '作成'
$data[22,0] = 'EditOperationName'
$data[22,1] = 'Edit'
$data[22,2] = '編集'
$data[23,0] = 'DeleteOperationName'
$data[23,1] = 'Delete'
$data[23,2] = '削除'
$data[24,0] = 'AddRemoveRobotsOperationName'
$data[24,1] = 'Add Or Remove Robots'
$data[24,2] = 'ロボットの追加・削除'
$data[25,0] = 'StoppedExecution'
$data[25,1] = 'The execution has been stopped.'
$data[25,2] = '実行が停止されました。'
$data[26,0] = ''
$data[26,1] = ''
$data[26,2] = ''
$data[27,0] = 'ChooseEntityMessage'
$data[27,1] = 'Choose entity to operate.'
$data[27,2] = 'エンティティを選択してください。'
$data[28,0] = 'ChooseEntityTitle'
$data[28,1] = 'Entity'
$data[28,2] = 'エンティティ'
$data[29,0] = 'ChooseOperationMessage'
$data[29,1] = 'Choose operation to perform.'
$data[29,2] = '操作を選択してください。'
$data[30,0] = 'ChooseOperationTitle'
$data[30,1] = 'Operation'
$data[30,2] = '操作'
$data[31,0] = 'UnsupportedEntity'
$data[31,1] = 'Unsupported entity specified.'
$data[31,2] = '選択されたエンティティが対応できません。'
$data[32,0] = 'UnsupportedOperation'
$data[32,1] = 'Unsupported operation.'
$data[32,2] = '選択された操作が対応できません。'
$data[33,0] = 'OperationDefaultResult'
$data[33,1] = 'Success'
$data[33,2] = '成功'
$data[34,0] = 'TokenNotRetrieved'
$data[34,1] = 'Unable to get access token due to failed authentication.'
$data[34,2] = 'API認証のためのトークンが取得できませんでした。'
$data[35,0] = 'ConfirmNumerousRequests'
$data[35,1] = 'The selected operation will make a large number of HTTP requests and might have a negative impact on infrastructure. Continue the processing?'
$data[35,2] = '選択された操作が多数のAPIリクエストを行う予想します。よろしいですか。'
$data[36,0] = 'NumerousRequestsNotConfirmed'
$data[36,1] = 'User did not confirm proceeding with numerous requests.'
$data[36,2] = 'ユーザーが多数のAPIリクエストを行うこと承認しませんでした。'
$data[37,0] = 'NumerousRequestsConfirmed'
$data[37,1] = 'User confirmed proceeding with numerous requests.'
$data[37,2] = 'ユーザーが多数のAPIリクエストを行うこと承認しました。'
$data[38,0] = 'ChooseOUMessage'
$data[38,1] = 'Choose Organization Unit'
$data[38,2] = '組織単位を選択してください。'
$data[39,0] = 'ChooseOUTitle'
$data[39,1] = 'Organization Unit'
$data[39,2] = '組織単位'
$data[40,0] = 'ProcessEntityFailure'
$data[40,1] = 'Failed to process entity. Request status: {0} / Response: {1}.'
$data[40,2] = '処理が失敗しました。リクエストステータス：{0} / レスポンス：{1}。'
$data[41,0] = 'ProcessedEntity'
$data[41,1] = 'Already processed.'
$data[41,2] = '処理済。'
$data[42,0] = 'IDInvalidOrNotSpecified'
$data[42,1] = 'ID invalid or not specified.'
$data[42,2] = 'IDが無効か指定されませんでした。'
$data[43,0] = 'NameNotSpecified'
$data[43,1] = 'Name not specified.'
$data[43,2] = '名前が指定されませんでした。'
$data[44,0] = 'IDAndNameDoNotMatch'
$data[44,1] = 'The specified ID and Name do not match.'
$data[44,2] = '指定された名前とIDが一致しません。'
$data[45,0] = 'TypeNotSpecified'
$data[45,1] = 'Type not specified.'
$data[45,2] = 'タイプが指定されませんでした。'
$data[46,0] = ''
$data[46,1] = ''
$data[46,2] = ''
$data[47,0] = 'UnsupportedAssetType'
$data[47,1] = 'The specified asset type is not supported.'
$data[47,2] = '指定されたアセットタイプが対応できません。'
$data[48,0] = 'AssetNotFound'
$data[48,1] = 'Asset not found.'
$data[48,2] = 'アセットが見つかりませんでした。'
$data[49,0] = 'AssetIDInvalidOrNotSpecified'
$data[49,1] = 'Asset ID invalid or not specified.'
$data[49,2] = 'アセットIDが無効か指定されませんでした。'
$data[50,0] = 'AssetNameNotSpecified'
$data[50,1] = 'Asset name not specified.'
$data[50,2] = 'アセット名が指定されませんでした。'
$data[51,0] = 'AssetIDAndNameDoNotMatch'
$data[51,1] = 'The specified Asset ID and Asset name do not match.'
$data[51,2] = '指定されたアセット名とアセットIDが一致しません。'
$data[52,0] = ''
$data[52,1] = ''
$data[52,2] = ''
$data[53,0] = 'UserNotFound'
$data[53,1] = 'User not  found.'
$data[53,2] = 'ユーザーが見つかりませんでした。'
$data[54,0] = 'StatusNotSupported'
$data[54,1] = 'Status not supported.'
$data[54,2] = '指定されたステータスが対応できません。'
$data[55,0] = 'UsernameNotSpecified'
$data[55,1] = 'Username not specified.'
$data[55,2] = 'ユーザー名が指定されませんでした。'
$data[56,0] = 'IDAndUsernameDoNotMatch'
$data[56,1] = 'The specified ID and Username do not match.'
$data[56,2] = '指定されたIDとユーザー名が一致しません。'
$data[57,0] = 'SurnameNotSpecified'
$data[57,1] = 'Surname not specified.'
$data[57,2] = '性が指定されませんでした。'
$data[58,0] = 'EmailNotSpecified'
$data[58,1] = 'Email not specified.'
$data[58,2] = 'メールアドレスが指定されませんでした。'
$data[59,0] = ''
$data[59,1] = ''
$data[59,2] = ''
$data[60,0] = 'ProcessedRobot'
$data[60,1] = 'Robot already processed. Robot name: {0} / ID: {1}.'
$data[60,2] = 'ロボット処理済。ロボット名：{0} / ID：{1}。'
$data[61,0] = 'RobotNotFound'
$data[61,1] = 'Robot not found.'
$data[61,2] = 'ロボットが見つかりませんでした。'
$data[62,0] = 'HostingTypeNotSpecified'
$data[62,1] = 'Hosting Type not specified.'
$data[62,2] = 'ホスティングの種類が指定されませんでした。'
$data[63,0] = 'RobotNameNotSpecified'
$data[63,1] = 'Robot Name not specified.'
$data[63,2] = 'ロボット名が指定されませんでした。'
$data[64,0] = 'RobotTypeNotSpecified'
$data[64,1] = 'Robot Type not specified.'
$data[64,2] = 'ロボットの種類が指定されませんでした。'
$data[65,0] = 'RobotIDInvalidOrNotSpecified'
$data[65,1] = 'Robot ID invalid or not specified.'
$data[65,2] = 'ロボットIDが無効か指定されませんでした。'
$data[66,0] = 'RobotIDAndNameDoNotMatch'
$data[66,1] = 'The specified Robot ID and name do not match.'
$data[66,2] = '指定されたロボット名とロボットIDが一致しません。'
$data[67,0] = 'NamedRobotNotFound'
$data[67,1] = 'Robot named {0} not found.'
$data[67,2] = '{0}というロボットが見つかりませんでした。'
$data[68,0] = ''
$data[68,1] = ''
$data[68,2] = ''
$data[69,0] = 'MachineNotFound'
$data[69,1] = 'Machine not found.'
$data[69,2] = 'マシンが見つかりませんでした。'
$data[70,0] = 'GetSingleMachineFailure'
$data[70,1] = 'Failed to get machine with Id: {0}. Request status: {1} / Response: {2}.'
$data[70,2] = 'ID{0}マシン取得が失敗しました。リクエストステータス：{1} / レスポンス：{2}。'
$data[71,0] = 'MachineNameNotSpecified'
$data[71,1] = 'Machine Name not specified.'
$data[71,2] = 'マシン名が指定されませんでした。'
$data[72,0] = ''
$data[72,1] = ''
$data[72,2] = ''
$data[73,0] = 'OUNotFound'
$data[73,1] = 'No Organization Unit found.'
$data[73,2] = '組織単位が見つかりませんでした。'
$data[74,0] = 'OUIDInvalidOrNotSpecified'
$data[74,1] = 'Organization Unit ID invalid or not specified.'
$data[74,2] = '組織単位IDが無効か指定されませんでした。'
$data[75,0] = 'OUNameNotSpecified'
$data[75,1] = 'Organization Unit name not specified.'
$data[75,2] = '組織単位名が指定されませんでした。'
$data[76,0] = 'OUIDAndNameDoNotMatch'
$data[76,1] = 'The specified Organization Unit ID and Organization Unit name do not match.'
$data[76,2] = '指定された組織単位名と組織単位IDが一致しません。'
$data[77,0] = ''
$data[77,1] = ''
$data[77,2] = ''
$data[78,0] = 'EnvironmentNotFound'
$data[78,1] = 'Environment not found.'
$data[78,2] = 'ロボットグループが見つかりませんでした。'
$data[79,0] = 'RobotEnvironmentAssociationFailure'
$data[79,1] = 'Failed to add robot to environment.　Request status: {0} / Response: {1}.'
$data[79,2] = 'ロボットをロボットグループに追加することができませんでした。リクエストステータス：{0} / レスポンス：{1}。'
$data[80,0] = 'EnvironmentIDInvalidOrNotSpecified'
$data[80,1] = 'Environment ID invalid or not specified.'
$data[80,2] = 'ロボットグループIDが無効か指定されませんでした。'
$data[81,0] = 'EnvironmentNameNotSpecified'
$data[81,1] = 'Environment name not specified.'
$data[81,2] = 'ロボットグループ名が指定されませんでした。'
$data[82,0] = 'EnvironmentIDAndNameDoNotMatch'
$data[82,1] = 'The specified Environment ID and Environment name do not match.'
$data[82,2] = '指定されたロボットグループ名とロボットグループIDが一致しません。'

$ws.Range("A1:C83").Value = $data

# Fix up cell formatting (wrap-text style) for the rows whose position
# within the reordered block does not already carry the right style.
$ws.Range("B71").WrapText = $true
$ws.Range("C71").WrapText = $true
$ws.Range("B76").WrapText = $true
$ws.Range("C76").WrapText = $true
$ws.Range("B78").WrapText = $true
$ws.Range("C78").WrapText = $true
$ws.Range("B79").WrapText = $true
$ws.Range("C79").WrapText = $true
$ws.Range("B82").WrapText = $true
$ws.Range("B80").ClearFormats()

# Resize the table to include the two newly-inserted rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C85"))

# Update the active cell selection to match the saved view state.
[void]$ws.Range("A2").Select()

Write-Output "done"
